$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the whole XML of a paragraph (found via its Range) with a
# modified version, built with a simple string substitution, so that we can
# make precise/structural OOXML edits while preserving all of the
# paragraph's other attributes (paraId, pPr, rsids, ...).
# ---------------------------------------------------------------------------
function Set-ParagraphXml($paraRange, [string]$oldFragment, [string]$newFragment) {
    $pxml = $paraRange.WordOpenXML
    $pStart = $pxml.IndexOf("<w:p ")
    if ($pStart -lt 0) { $pStart = $pxml.IndexOf("<w:p>") }
    $pEndTag = "</w:p>"
    $pEnd = $pxml.IndexOf($pEndTag, $pStart) + $pEndTag.Length
    $fullP = $pxml.Substring($pStart, $pEnd - $pStart)

    $newP = $fullP.Replace($oldFragment, $newFragment)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:body>' + $newP + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$paraRange.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) The first three inline pictures just gain <w:noProof/> in their run's
#    rPr (existing b/bCs is kept as-is).
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(1).Range.Font.NoProofing = 1
$d.InlineShapes.Item(2).Range.Font.NoProofing = 1
$d.InlineShapes.Item(3).Range.Font.NoProofing = 1

# ---------------------------------------------------------------------------
# 2) The floating (anchored) picture right after "Program 2" also gains
#    <w:noProof/> (keeping b/bCs). It is not part of InlineShapes, so locate
#    its paragraph via the preceding "Program 2" text.
# ---------------------------------------------------------------------------
$prg2 = $d.Content
$prg2.Find.ClearFormatting()
$prg2.Find.Text = "Program 2"
[void]$prg2.Find.Execute()

$paraCount = $d.Paragraphs.Count
$prg2ParaIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $prg2.Start -and $p.Range.End -ge $prg2.End) {
        $prg2ParaIdx = $i
    }
}
$floatPara = $d.Paragraphs.Item($prg2ParaIdx + 1)
$floatPara.Range.Font.NoProofing = 1

# ---------------------------------------------------------------------------
# 3) The second "Sample input and output" run gets a <w:lastRenderedPageBreak/>
#    inserted right before its <w:t>.
# ---------------------------------------------------------------------------
$sample = $d.Content
$sample.Find.ClearFormatting()
$sample.Find.Text = "Sample input and output"
[void]$sample.Find.Execute()
$sample.Collapse(0)
[void]$sample.Find.Execute()

$paraCount = $d.Paragraphs.Count
$sampleParaIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $sample.Start -and $p.Range.End -ge $sample.End) {
        $sampleParaIdx = $i
    }
}
$samplePara = $d.Paragraphs.Item($sampleParaIdx)
Set-ParagraphXml $samplePara.Range "<w:t>Sample input and output</w:t>" "<w:lastRenderedPageBreak/><w:t>Sample input and output</w:t>"

# ---------------------------------------------------------------------------
# 4) The fourth inline picture (the one that now directly follows the run
#    edited above) has its run rPr fully replaced: b/bCs/u all removed and
#    replaced by a single <w:noProof/>.
# ---------------------------------------------------------------------------
$lastShapeRng = $d.InlineShapes.Item(4).Range
$paraCount = $d.Paragraphs.Count
$lastShapeParaIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $lastShapeRng.Start -and $p.Range.End -ge $lastShapeRng.End) {
        $lastShapeParaIdx = $i
    }
}
$lastShapePara = $d.Paragraphs.Item($lastShapeParaIdx)
Set-ParagraphXml $lastShapePara.Range "<w:rPr><w:b/><w:bCs/><w:u w:val=`"single`"/></w:rPr><w:drawing" "<w:rPr><w:noProof/></w:rPr><w:drawing"
